$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete higher-numbered row first so row numbers below stay stable while we work
# Row 28 = "SC 92" -> delete entirely
$ws.Rows.Item(28).Delete()
# Row 26 = "RM 232" -> delete entirely
$ws.Rows.Item(26).Delete()

# After both deletions, rows shifted up by 2:
# Row 26 = SC 5   -> E26 becomes -5
# Row 27 = SC 101 -> E27 becomes empty
# Row 28 = SC 105 (unchanged)
# Row 29 = SC 119 (unchanged)
# Row 30 = SC 120 -> E30 becomes -5.7
# Row 31 = SC 132 (unchanged)
# Row 32 = SC 193 -> E32 becomes empty
# Row 33 = SC 232 (unchanged)

$ws.Range("E26").Value = -5
$ws.Range("E27").Value = ""
$ws.Range("E30").Value = -5.7
$ws.Range("E32").Value = ""
